# The deck's design was switched from the custom "Integral" theme back to
# the plain default "Office Theme" palette (dk1/lt1/dk2/lt2/accent1-6/
# hlink/folHlink). Re-apply each of the twelve theme colour slots through
# the Master's ColorScheme, which is what PowerPoint COM automation uses
# to push new RGB values into the presentation's theme part.

$p = $ppt.ActivePresentation
$master = $p.SlideMaster
$colorScheme = $master.ColorScheme

function Set-ThemeColor($index, $r, $g, $b) {
    $rgb = $r + ($g * 256) + ($b * 65536)
    $colorScheme.Colors($index).RGB = $rgb
}

# 1  dk1      000000
Set-ThemeColor 1  0x00 0x00 0x00
# 2  lt1      FFFFFF
Set-ThemeColor 2  0xFF 0xFF 0xFF
# 3  dk2      44546A
Set-ThemeColor 3  0x44 0x54 0x6A
# 4  lt2      E7E6E6
Set-ThemeColor 4  0xE7 0xE6 0xE6
# 5  accent1  5B9BD5
Set-ThemeColor 5  0x5B 0x9B 0xD5
# 6  accent2  ED7D31
Set-ThemeColor 6  0xED 0x7D 0x31
# 7  accent3  A5A5A5
Set-ThemeColor 7  0xA5 0xA5 0xA5
# 8  accent4  FFC000
Set-ThemeColor 8  0xFF 0xC0 0x00
# 9  accent5  4472C4
Set-ThemeColor 9  0x44 0x72 0xC4
# 10 accent6  70AD47
Set-ThemeColor 10 0x70 0xAD 0x47
# 11 hlink    0563C1
Set-ThemeColor 11 0x05 0x63 0xC1
# 12 folHlink 954F72
Set-ThemeColor 12 0x95 0x4F 0x72
